# Modification to Benders' main cycle.
# Update the Scenarios sheet inputs: numScenarios (B1), and the Upward/Downward
# multipliers in column C (C1, C2, C3) all become 1. Dependent formulas across
# the DownActivation/UpActivation sheets (which reference Scenarios!$C$2) will
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1

# Move the active selection to R6, matching the author's last cursor position.
$ws.Range("R6").Select()
